# "Trabajando en logica de finalizado ruta"
#
# 1) clientes: update I3 (ULTIMO RETIRO) and K3 (PROXIMO RETIRO) for the
#    Isaias Beroiza Mora row.
# 2) rutas_bd: append three new route-log rows (10-12).
# 3) rutas_registros: fill in the two still-open "ruta ejemplo"/"ruta test"
#    rows (16-17) and append the closing rows (18-19) plus a trailing blank
#    row (20).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# clientes
# ---------------------------------------------------------------------
$clientes = $wb.Worksheets.Item("clientes")

$clientes.Range("I3").NumberFormat = "@"
$clientes.Range("I3").Value = "20240808"

$clientes.Range("K3").NumberFormat = "@"
$clientes.Range("K3").Value = "2024-12-08"

# ---------------------------------------------------------------------
# rutas_bd
# ---------------------------------------------------------------------
$rutasBd = $wb.Worksheets.Item("rutas_bd")

function Set-RutasBdRow($row, $fecha, $situacion, $otro) {
    $rutasBd.Range("A$row").NumberFormat = "@"
    $rutasBd.Range("A$row").Value = $fecha
    $rutasBd.Range("B$row").Value = 1
    $rutasBd.Range("C$row").Value = "16.742.249-7"
    $rutasBd.Range("D$row").Value = "Isaias Beroiza Mora"
    $rutasBd.Range("E$row").Value = "colaco sn km3 parcela 9"
    $rutasBd.Range("F$row").Value = "Calbuco"
    $rutasBd.Range("G$row").NumberFormat = "@"
    $rutasBd.Range("G$row").Value = "88809703"
    $rutasBd.Range("H$row").Value = "por buscar"
    $rutasBd.Range("I$row").Value = "ok"
    $rutasBd.Range("J$row").Value = $situacion
    $rutasBd.Range("K$row").Value = $otro
}

Set-RutasBdRow 10 "20240829" "REALIZADO" "3c3l1ba ok"
Set-RutasBdRow 11 "20240828" "REALIZADO" "3c3l1ba ok"
Set-RutasBdRow 12 "20240808" "POSPUESTO" "DEUDA"

# ---------------------------------------------------------------------
# rutas_registros
# ---------------------------------------------------------------------
$rutasReg = $wb.Worksheets.Item("rutas_registros")

# Columns B, D and E carry a sheet-level column style; new cells written
# into them need to be nudged back to the plain/default look so they
# don't pick up that column formatting.
function Set-PlainStyle($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
    $range.VerticalAlignment = -4108
}

# Row 16: close out "ruta ejemplo" (20240829)
$rutasReg.Range("A16").NumberFormat = "@"
$rutasReg.Range("A16").Value = "20240829"
$rutasReg.Range("B16").Value = "ruta ejemplo"
Set-PlainStyle $rutasReg.Range("B16")
$rutasReg.Range("C16").Value = 1
$rutasReg.Range("F16").Value = "RUTA FINALIZADA"

# Row 17: close out "ruta test" (20240828)
$rutasReg.Range("A17").NumberFormat = "@"
$rutasReg.Range("A17").Value = "20240828"
$rutasReg.Range("B17").Value = "ruta test"
Set-PlainStyle $rutasReg.Range("B17")
$rutasReg.Range("C17").Value = 1
$rutasReg.Range("F17").Value = "RUTA FINALIZADA"

# Row 18: another "ruta test" closure (20240807)
$rutasReg.Range("A18").NumberFormat = "@"
$rutasReg.Range("A18").Value = "20240807"
$rutasReg.Range("B18").Value = "ruta test"
Set-PlainStyle $rutasReg.Range("B18")
$rutasReg.Range("F18").Value = "RUTA FINALIZADA"

# Row 19: "Ruta ejemplo terminacion" closure, with a timestamped note
$rutasReg.Range("A19").NumberFormat = "@"
$rutasReg.Range("A19").Value = "20240808"
$rutasReg.Range("B19").Value = "Ruta ejemplo terminacion"
Set-PlainStyle $rutasReg.Range("B19")
$rutasReg.Range("D19").Value = 1
$rutasReg.Range("F19").Value = "2024-08-06T15:51:36.376168 RUTA FINALIZADA"

# Row 20: trailing blank row (keeps the used range extending to F20)
$rutasReg.Range("F20").Value = ""
